$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H130").Value = 47224
$ws.Range("J130").Value = 47224
$ws.Range("L130").Value = 47224
$ws.Range("N130").Value = -57264
$ws.Range("H137").Value = 3671.3137
$ws.Range("I137").Value = 1212.5676
$ws.Range("J137").Value = 10169.429
$ws.Range("K137").Value = 3637.7028
$ws.Range("L137").Value = 30508.287
$ws.Range("M137").Value = -1087.7028
$ws.Range("N137").Value = -35608.287

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 49111.715
$ws.Range("J80").Value = 49111.715
$ws.Range("L80").Value = 49111.715
$ws.Range("N80").Value = -51107.715
$ws.Range("H83").Value = 49111.715
$ws.Range("J83").Value = 49111.715
$ws.Range("L83").Value = 147335.145
$ws.Range("N83").Value = -157319.145
$ws.Range("H137").Value = 44645
$ws.Range("J137").Value = 44645
$ws.Range("L137").Value = 44645
$ws.Range("N137").Value = -54845
$ws.Range("H138").Value = 32711
$ws.Range("J138").Value = 32711
$ws.Range("L138").Value = 32711
$ws.Range("N138").Value = -42991

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 39095
$ws.Range("J59").Value = 39095
$ws.Range("L59").Value = 39095
$ws.Range("N59").Value = -40789
$ws.Range("H124").Value = 49983.332
$ws.Range("J124").Value = 49983.332
$ws.Range("L124").Value = 49983.332
$ws.Range("N124").Value = -59803.332
$ws.Range("H133").Value = 57389
$ws.Range("J133").Value = 57389
$ws.Range("L133").Value = 57389
$ws.Range("N133").Value = -67509
$ws.Range("H139").Value = 63299.668
$ws.Range("J139").Value = 63299.668
$ws.Range("L139").Value = 63299.668
$ws.Range("N139").Value = -73579.66800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3049.5579
$ws.Range("I31").Value = 1109.5927
$ws.Range("J31").Value = 3819.8381
$ws.Range("K31").Value = 1109.5927
$ws.Range("L31").Value = 3819.8381
$ws.Range("M31").Value = -814.5926999999999
$ws.Range("N31").Value = -4409.8381
$ws.Range("H34").Value = 3049.5579
$ws.Range("I34").Value = 1109.5927
$ws.Range("J34").Value = 3819.8381
$ws.Range("K34").Value = 1109.5927
$ws.Range("L34").Value = 3819.8381
$ws.Range("M34").Value = -907.5926999999999
$ws.Range("N34").Value = -4223.8381
$ws.Range("H58").Value = 1259.4681
$ws.Range("I58").Value = 869.11536
$ws.Range("J58").Value = 1742.762
$ws.Range("K58").Value = 869.11536
$ws.Range("L58").Value = 1742.762
$ws.Range("M58").Value = -666.11536
$ws.Range("N58").Value = -2148.762
$ws.Range("H81").Value = 48311
$ws.Range("J81").Value = 48311
$ws.Range("L81").Value = 48311
$ws.Range("N81").Value = -50307
$ws.Range("H82").Value = 26073
$ws.Range("I82").Value = 3000
$ws.Range("K82").Value = 3000
$ws.Range("M82").Value = -2639
$ws.Range("H84").Value = 48311
$ws.Range("J84").Value = 48311
$ws.Range("L84").Value = 144933
$ws.Range("N84").Value = -154917
$ws.Range("H85").Value = 26073
$ws.Range("I85").Value = 3000
$ws.Range("K85").Value = 3000
$ws.Range("M85").Value = -1752
$ws.Range("H88").Value = 21778.334
$ws.Range("J88").Value = 21778.334
$ws.Range("L88").Value = 21778.334
$ws.Range("N88").Value = -22590.334
$ws.Range("H91").Value = 21778.334
$ws.Range("J91").Value = 21778.334
$ws.Range("L91").Value = 21778.334
$ws.Range("N91").Value = -24586.334
$ws.Range("H100").Value = 47139.25
$ws.Range("J100").Value = 47139.25
$ws.Range("L100").Value = 47139.25
$ws.Range("N100").Value = -49303.25
$ws.Range("H136").Value = 1259.4681
$ws.Range("I136").Value = 869.11536
$ws.Range("J136").Value = 1742.762
$ws.Range("K136").Value = 2607.34608
$ws.Range("L136").Value = 5228.286
$ws.Range("M136").Value = -57.34608000000026
$ws.Range("N136").Value = -10328.286
$ws.Range("H137").Value = 70604.22
$ws.Range("J137").Value = 70604.22
$ws.Range("L137").Value = 70604.22
$ws.Range("N137").Value = -80804.22

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 20779.309
$ws.Range("J74").Value = 20779.309
$ws.Range("L74").Value = 20779.309
$ws.Range("N74").Value = -22651.309
$ws.Range("H77").Value = 20779.309
$ws.Range("J77").Value = 20779.309
$ws.Range("L77").Value = 62337.927
$ws.Range("N77").Value = -71697.927
$ws.Range("H80").Value = 159872.03
$ws.Range("I80").Value = 281467.62
$ws.Range("J80").Value = 3534.8572
$ws.Range("K80").Value = 281467.62
$ws.Range("L80").Value = 3534.8572
$ws.Range("M80").Value = -280469.62
$ws.Range("N80").Value = -5530.8572
$ws.Range("H83").Value = 159872.03
$ws.Range("I83").Value = 281467.62
$ws.Range("J83").Value = 3534.8572
$ws.Range("K83").Value = 1407338.1
$ws.Range("L83").Value = 17674.286
$ws.Range("M83").Value = -1402346.1
$ws.Range("N83").Value = -27658.286
$ws.Range("H137").Value = 19346.666
$ws.Range("J137").Value = 19346.666
$ws.Range("L137").Value = 19346.666
$ws.Range("N137").Value = -29546.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H88").Value = 37383.75
$ws.Range("J88").Value = 43178.332
$ws.Range("L88").Value = 43178.332
$ws.Range("N88").Value = -44034.332
$ws.Range("H91").Value = 37383.75
$ws.Range("J91").Value = 43178.332
$ws.Range("L91").Value = 43178.332
$ws.Range("N91").Value = -46142.332
$ws.Range("H92").Value = 38499
$ws.Range("J92").Value = 38499
$ws.Range("L92").Value = 38499
$ws.Range("N92").Value = -43491
$ws.Range("H96").Value = 28846.5
$ws.Range("J96").Value = 28846.5
$ws.Range("L96").Value = 28846.5
$ws.Range("N96").Value = -34338.5
$ws.Range("H99").Value = 26022.9
$ws.Range("I99").Value = 18045.8
$ws.Range("J99").Value = 34000
$ws.Range("K99").Value = 18045.8
$ws.Range("L99").Value = 34000
$ws.Range("M99").Value = -15050.8
$ws.Range("N99").Value = -39990
$ws.Range("H102").Value = 48553
$ws.Range("J102").Value = 48553
$ws.Range("L102").Value = 48553
$ws.Range("N102").Value = -55043
$ws.Range("H137").Value = 29100
$ws.Range("J137").Value = 29100
$ws.Range("L137").Value = 29100
$ws.Range("N137").Value = -39300
$ws.Range("H139").Value = 64671.285
$ws.Range("J139").Value = 48783.168
$ws.Range("L139").Value = 48783.168
$ws.Range("N139").Value = -59063.168

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 24020
$ws.Range("J92").Value = 24020
$ws.Range("L92").Value = 24020
$ws.Range("N92").Value = -29012
$ws.Range("H93").Value = 36114.832
$ws.Range("J93").Value = 36114.832
$ws.Range("L93").Value = 36114.832
$ws.Range("N93").Value = -41106.832
$ws.Range("H99").Value = 37286.4
$ws.Range("J99").Value = 36375
$ws.Range("L99").Value = 36375
$ws.Range("N99").Value = -42365
$ws.Range("H102").Value = 41329
$ws.Range("J102").Value = 41329
$ws.Range("L102").Value = 41329
$ws.Range("N102").Value = -47819
$ws.Range("H106").Value = 33998.855
$ws.Range("J106").Value = 33998.855
$ws.Range("L106").Value = 33998.855
$ws.Range("N106").Value = -36522.855
$ws.Range("H133").Value = 77114.25
$ws.Range("J133").Value = 77114.25
$ws.Range("L133").Value = 77114.25
$ws.Range("N133").Value = -87234.25
$ws.Range("H136").Value = 18834.91
$ws.Range("I136").Value = 27622.945
$ws.Range("J136").Value = 1721.3684
$ws.Range("K136").Value = 82868.83499999999
$ws.Range("L136").Value = 5164.1052
$ws.Range("M136").Value = -80318.83499999999
$ws.Range("N136").Value = -10264.1052
$ws.Range("H137").Value = 43818.848
$ws.Range("J137").Value = 43818.848
$ws.Range("L137").Value = 43818.848
$ws.Range("N137").Value = -54018.848
